$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Hora (G) columns keep their original text formatting
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Updated prices for Wed Dec 14 13:24:28 UTC 2022 run
$ws.Range("D2").Value = "270.54"
$ws.Range("D3").Value = "22.94"
$ws.Range("D4").Value = "6.394"
$ws.Range("D5").Value = "0.06262"
$ws.Range("D6").Value = "3.647"
$ws.Range("D7").Value = "6.694"
$ws.Range("D8").Value = "1.375"
$ws.Range("D9").Value = "0.8353"
$ws.Range("D10").Value = "0.01380"
$ws.Range("D11").Value = "0.1630"
$ws.Range("D12").Value = "0.08387"
$ws.Range("D13").Value = "0.03387"
$ws.Range("D14").Value = "0.03098"
$ws.Range("D15").Value = "0.1264"
$ws.Range("D16").Value = "0.09333"
$ws.Range("D17").Value = "3.868"
$ws.Range("D18").Value = "0.001640"
$ws.Range("D19").Value = "0.04801"
$ws.Range("D20").Value = "0.006300"
$ws.Range("D22").Value = "0.003329"
$ws.Range("D24").Value = "3.736"
$ws.Range("D25").Value = "2.369"
$ws.Range("D26").Value = "0.3404"
$ws.Range("D41").Value = "0.006923"
$ws.Range("D43").Value = "0.003419"
$ws.Range("D44").Value = "0.01229"
$ws.Range("D45").Value = "0.00006281"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("D47").Value = "0.9000"
$ws.Range("D48").Value = "0.06228"
$ws.Range("D49").Value = "0.00002298"

# Hora column bumped from 12 to 13 for every data row
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 7).Value = "13"
}
